# TC01_Trials_Filter_Ethnicity-HispLatino.xlsx — add the Neo4j/web query text
# to the "startup" sheet (A2), resize row 2 to fit it, and update the
# worksheet selection to match the post-edit authoring session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A2 previously held no value (just the wrap-text style). Populate it with the
# Cypher/Neo4j query used to build this test-case's data extract.
$ws.Range("A2").Value = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.ethnicity IN [''HISPANIC_OR_LATINO''] RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'

# The long wrapped query text needs a taller row to display.
$ws.Rows.Item(2).RowHeight = 87

# Update the live selection: the sheet no longer scrolls right to column B
# (topLeftCell reset to default) and the selected range grows to B2:B6.
$ws.Range("B2:B6").Select() | Out-Null
